$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the text of B15: append the new clause describing that the
#    sentence-transformer embedding was added to the elasticSearch db.
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "huggingface: init, save, load sentence transformer, added as embedding to elasticSearch db"

# ---------------------------------------------------------------------------
# 2. Give the rows for 2023-08-11 .. 2023-08-13 (rows 13-15) the same boxed
#    "single day" border treatment already used for the other day blocks
#    further up in the sheet, by copying the cell formats from rows that
#    already carry the desired look (this reuses the existing style
#    definitions instead of creating new duplicate ones).
# ---------------------------------------------------------------------------

# Row 13 becomes the top edge of the box (like row 6 / row 8 for column A,
# and row 6 for column B).
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

# Row 14 reuses the full single-row box style (like rows 9-12).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

# Row 15 becomes the bottom edge of the box (like row 5 for column A,
# and row 7 for column B).
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Update row heights to match the new boxed layout / wrapped text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 52
$ws.Rows.Item(14).RowHeight = 35
$ws.Rows.Item(15).RowHeight = 35

# ---------------------------------------------------------------------------
# 4. Update the active selection to B19 (where the user ended up after the
#    edit), keeping the same view.
# ---------------------------------------------------------------------------
$ws.Range("B19").Select() | Out-Null
